$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both need F2, F4, F5 updated.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 135
    $ws.Range("F4").Value = 88
    $ws.Range("F5").Value = 28
}
